$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.791.37"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.571.69"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "563.12"
$ws.Range("E5").Value = "  +1.65%  "
$ws.Range("D6").Value = "142.86"
$ws.Range("E6").Value = "  -1.40%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "2.577.98"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "6.66"
$ws.Range("E10").Value = "  -2.36%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("E12").Value = "  +8.36%  "
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("D14").Value = "3.024.85"
$ws.Range("E14").Value = "  -0.86%  "
$ws.Range("D15").Value = "58.877.72"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("E16").Value = "  +5.41%  "
$ws.Range("E17").Value = "  +2.95%  "
$ws.Range("D18").Value = "2.563.69"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "334.53"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "63.98"
$ws.Range("E24").Value = "  -3.98%  "
$ws.Range("E25").Value = "  +5.67%  "
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "0.0₃0774"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -0.67%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "157.94"
$ws.Range("E32").Value = "  +2.32%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "6.02"
$ws.Range("E33").Value = "  +0.56%  "
$ws.Range("D34").Value = "18.93"
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "0.870"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "0.874"
$ws.Range("E37").Value = "  +5.61%  "
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "36.76"
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("E40").Value = "  +1.61%  "
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "289.42"
$ws.Range("E42").Value = "  +1.90%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "0.0968"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").Value = "0.593"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("D46").Value = "10.62"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("D49").Value = "123.91"
$ws.Range("E49").Value = "  +7.81%  "
$ws.Range("E50").Value = "  +1.35%  "
$ws.Range("D51").Value = "18.44"
$ws.Range("E51").Value = "  +2.15%  "
